# run_settings.xlsx edit
# - bump D2:D5 on Sheet1 from 0 -> 0.1
# - update sheetView selection on Sheet1
# - add new worksheet "Stepwise_reduction_waste" after Sheet1, populate it,
#   make it the active/selected sheet
# - absPath folder rename + revisionPtr/window-geometry are session/host
#   artifacts outside the reach of the exposed object model; skipped.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: dev column 0 -> 0.1 ------------------------------------------
$ws1.Range("D2").Value = 0.1
$ws1.Range("D3").Value = 0.1
$ws1.Range("D4").Value = 0.1
$ws1.Range("D5").Value = 0.1

# --- add the new sheet right after Sheet1 ----------------------------------
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "Stepwise_reduction_waste"

# Column widths matching Sheet1's layout
$ws2.Columns.Item(2).ColumnWidth = 10.28515625
$ws2.Columns.Item(3).ColumnWidth = 10.28515625
$ws2.Columns.Item(4).ColumnWidth = 9.42578125
$ws2.Columns.Item(5).ColumnWidth = 14
$ws2.Columns.Item(6).ColumnWidth = 22.42578125

# Run labels (column A) written first & top-to-bottom so the shared-string
# table grows in the same order Excel produced it: run, run 1.. run 4 reuse
# existing strings, run 5..run 16 are newly appended (indices 15-26).
$runLabels = @("run","run 1","run 2","run 3","run 4","run 5","run 6","run 7","run 8","run 9","run 10","run 11","run 12","run 13","run 14","run 15","run 16")
for ($i = 0; $i -lt $runLabels.Length; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $runLabels[$i]
}

# Column B (n_days / 5)
for ($r = 1; $r -le 17; $r++) {
    if ($r -eq 1) {
        $ws2.Cells.Item($r, 2).Value = "n_days"
    } else {
        $ws2.Cells.Item($r, 2).Value = 5
    }
}

# Column C (n_persons / 4)
for ($r = 1; $r -le 17; $r++) {
    if ($r -eq 1) {
        $ws2.Cells.Item($r, 3).Value = "n_persons"
    } else {
        $ws2.Cells.Item($r, 3).Value = 4
    }
}

# Column D (dev / 0.1)
for ($r = 1; $r -le 17; $r++) {
    if ($r -eq 1) {
        $ws2.Cells.Item($r, 4).Value = "dev"
    } else {
        $ws2.Cells.Item($r, 4).Value = 0.1
    }
}

# Column E (optimize_over / Total_carbon)
$ws2.Cells.Item(1, 5).Value = "optimize_over"
for ($r = 2; $r -le 17; $r++) {
    $ws2.Cells.Item($r, 5).Value = "Total_carbon"
}

# Column F (DRVs / modelgezin_gemiddeld)
$ws2.Cells.Item(1, 6).Value = "DRVs"
for ($r = 2; $r -le 17; $r++) {
    $ws2.Cells.Item($r, 6).Value = "modelgezin_gemiddeld"
}

# Column G (tvar1 header + numeric series)
$ws2.Cells.Item(1, 7).Value = "tvar1"
$gValues = @(890, 810, 720, 630, 540, 450, 360, 270, 180, 90, 80, 70, 60, 50, 40, 39)
for ($i = 0; $i -lt $gValues.Length; $i++) {
    $ws2.Cells.Item($i + 2, 7).Value = $gValues[$i]
}

# Bold header row (reuses the existing bold style already used on Sheet1 row 1)
$ws2.Range("A1:G1").Font.Bold = $true

# --- selections / active sheet ---------------------------------------------
$null = $ws1.Range("D6").Select()
$ws2.Activate()
$null = $ws2.Range("G2").Select()
